# Apply updated crypto price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# get a leading "'" (quote-prefix) so they stay plain text, matching the
# original inlineStr cells.

$ws.Range("D2").Value = "25.811.72"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.639.13"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'215.51"
$ws.Range("D6").Value = "'0.5070"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "'0.2583"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").Value = "'0.06421"
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("D10").Value = "'20.39"
$ws.Range("E10").Value = "  +4.80%  "
$ws.Range("D11").Value = "'0.07782"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "'4.256"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "1.644.40"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").Value = "1.867.06"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").Value = "'0.5610"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").Value = "0.0₅7647"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "'63.33"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "25.838.07"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").Value = "'4.371"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("D21").Value = "'192.18"
$ws.Range("E21").Value = "  -1.37%  "
$ws.Range("D22").Value = "'9.912"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("D23").Value = "'6.143"
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "'1.791"
$ws.Range("E25").Value = "  -5.29%  "
$ws.Range("D26").Value = "'139.68"
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("D27").Value = "'0.1229"
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("D28").Value = "'6.813"
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("D29").Value = "'15.54"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").Value = "'1.244"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").Value = "'0.04938"
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("D32").Value = "'3.281"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").Value = "'3.248"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("D34").Value = "'1.569"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("D35").Value = "'2.386"
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("D36").Value = "'0.9032"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("E37").Value = "  +1.45%  "
$ws.Range("D38").Value = "'0.5568"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("D39").Value = "1.132.28"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").Value = "'0.9969"
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").Value = "'5.457"
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("D43").Value = "'0.8002"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").Value = "'98.91"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("D45").Value = "1.778.56"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").Value = "0.0₈112"
$ws.Range("E46").Value = "  -4.26%  "
$ws.Range("D47").Value = "'55.55"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("E48").Value = "  -3.75%  "
$ws.Range("D49").Value = "'7.783"
$ws.Range("E49").Value = "  +2.71%  "
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "  +0.51%  "
